$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(2, 2).Value = "701号直流"
$ws.Cells.Item(2, 3).Value = 45927.457337962966
$ws.Cells.Item(2, 4).Value = 45935.328888888886
$ws.Cells.Item(2, 5).Value = 188.91722222208045

$ws.Cells.Item(3, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(3, 2).Value = "502号直流"
$ws.Cells.Item(3, 3).Value = 45930.238043981481
$ws.Cells.Item(3, 4).Value = 45935.328888888886
$ws.Cells.Item(3, 5).Value = 122.18027777771931

$ws.Cells.Item(4, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(4, 2).Value = "112号直流"
$ws.Cells.Item(4, 3).Value = 45930.517060185186
$ws.Cells.Item(4, 4).Value = 45935.328888888886
$ws.Cells.Item(4, 5).Value = 115.48388888878981

$ws.Cells.Item(5, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(5, 2).Value = "111号直流"
$ws.Cells.Item(5, 3).Value = 45930.618518518517
$ws.Cells.Item(5, 4).Value = 45935.328888888886
$ws.Cells.Item(5, 5).Value = 113.04888888885034

$ws.Cells.Item(6, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(6, 2).Value = "201号直流"
$ws.Cells.Item(6, 3).Value = 45931.575543981482
$ws.Cells.Item(6, 4).Value = 45935.328888888886
$ws.Cells.Item(6, 5).Value = 90.080277777684387

$ws.Cells.Item(7, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(7, 2).Value = "603号直流"
$ws.Cells.Item(7, 3).Value = 45932.081099537034
$ws.Cells.Item(7, 4).Value = 45935.328888888886
$ws.Cells.Item(7, 5).Value = 77.946944444440305

$ws.Cells.Item(8, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(8, 2).Value = "406号直流"
$ws.Cells.Item(8, 3).Value = 45933.039143518516
$ws.Cells.Item(8, 4).Value = 45935.328888888886
$ws.Cells.Item(8, 5).Value = 54.953888888878282

$ws.Cells.Item(9, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(9, 2).Value = "904号直流"
$ws.Cells.Item(9, 3).Value = 45933.088784722226
$ws.Cells.Item(9, 4).Value = 45935.328888888886
$ws.Cells.Item(9, 5).Value = 53.762499999837019

$ws.Cells.Item(10, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(10, 2).Value = "504号直流"
$ws.Cells.Item(10, 3).Value = 45933.270925925928
$ws.Cells.Item(10, 4).Value = 45935.328888888886
$ws.Cells.Item(10, 5).Value = 49.391111110977363

$ws.Cells.Item(11, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(11, 2).Value = "103号直流"
$ws.Cells.Item(11, 3).Value = 45933.305023148147
$ws.Cells.Item(11, 4).Value = 45935.328888888886
$ws.Cells.Item(11, 5).Value = 48.57277777773561

$ws.Cells.Item(12, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(12, 2).Value = "002B号直流"
$ws.Cells.Item(12, 3).Value = 45933.517557870371
$ws.Cells.Item(12, 4).Value = 45935.328888888886
$ws.Cells.Item(12, 5).Value = 43.471944444347173

$ws.Cells.Item(13, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(13, 2).Value = "103号直流"
$ws.Cells.Item(13, 3).Value = 45934.067881944444
$ws.Cells.Item(13, 4).Value = 45935.328888888886
$ws.Cells.Item(13, 5).Value = 30.26416666660225

$ws.Cells.Item(14, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(14, 2).Value = "503号直流"
$ws.Cells.Item(14, 3).Value = 45934.068425925929
$ws.Cells.Item(14, 4).Value = 45935.328888888886
$ws.Cells.Item(14, 5).Value = 30.251111110963393

$ws.Cells.Item(15, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(15, 2).Value = "203号直流"
$ws.Cells.Item(15, 3).Value = 45934.232268518521
$ws.Cells.Item(15, 4).Value = 45935.328888888886
$ws.Cells.Item(15, 5).Value = 26.318888888752554

$ws.Cells.Item(16, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(16, 2).Value = "003B号直流"
$ws.Cells.Item(16, 3).Value = 45934.262986111113
$ws.Cells.Item(16, 4).Value = 45935.328888888886
$ws.Cells.Item(16, 5).Value = 25.581666666548699

$ws.Cells.Item(17, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(17, 2).Value = "505号直流"
$ws.Cells.Item(17, 3).Value = 45934.284421296295
$ws.Cells.Item(17, 4).Value = 45935.328888888886
$ws.Cells.Item(17, 5).Value = 25.067222222161945

$ws.Cells.Item(18, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(18, 2).Value = "106号直流"
$ws.Cells.Item(18, 3).Value = 45934.342083333337
$ws.Cells.Item(18, 4).Value = 45935.328888888886
$ws.Cells.Item(18, 5).Value = 23.683333333174232

$ws.Cells.Item(19, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(19, 2).Value = "702号直流"
$ws.Cells.Item(19, 3).Value = 45934.344675925924
$ws.Cells.Item(19, 4).Value = 45935.328888888886
$ws.Cells.Item(19, 5).Value = 23.621111111075152

$ws.Cells.Item(20, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(20, 2).Value = "006B号直流"
$ws.Cells.Item(20, 3).Value = 45934.406018518515
$ws.Cells.Item(20, 4).Value = 45935.328888888886
$ws.Cells.Item(20, 5).Value = 22.148888888885267

$ws.Cells.Item(21, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(21, 2).Value = "901号直流"
$ws.Cells.Item(21, 3).Value = 45934.527488425927
$ws.Cells.Item(21, 4).Value = 45935.328888888886
$ws.Cells.Item(21, 5).Value = 19.233611111005303

$ws.Cells.Item(22, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(22, 2).Value = "102号直流"
$ws.Cells.Item(22, 3).Value = 45934.532418981478
$ws.Cells.Item(22, 4).Value = 45935.328888888886
$ws.Cells.Item(22, 5).Value = 19.115277777775191

$ws.Cells.Item(23, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(23, 2).Value = "903号直流"
$ws.Cells.Item(23, 3).Value = 45934.542071759257
$ws.Cells.Item(23, 4).Value = 45935.328888888886
$ws.Cells.Item(23, 5).Value = 18.883611111086793

$ws.Cells.Item(24, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(24, 2).Value = "801号直流"
$ws.Cells.Item(24, 3).Value = 45934.554988425924
$ws.Cells.Item(24, 4).Value = 45935.328888888886
$ws.Cells.Item(24, 5).Value = 18.573611111089122

$ws.Cells.Item(25, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(25, 2).Value = "001B号直流"
$ws.Cells.Item(25, 3).Value = 45934.55914351852
$ws.Cells.Item(25, 4).Value = 45935.328888888886
$ws.Cells.Item(25, 5).Value = 18.473888888780493

$ws.Cells.Item(26, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(26, 2).Value = "905号直流"
$ws.Cells.Item(26, 3).Value = 45934.572893518518
$ws.Cells.Item(26, 4).Value = 45935.328888888886
$ws.Cells.Item(26, 5).Value = 18.143888888822403

$ws.Cells.Item(27, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(27, 2).Value = "102号直流"
$ws.Cells.Item(27, 3).Value = 45934.578287037039
$ws.Cells.Item(27, 4).Value = 45935.328888888886
$ws.Cells.Item(27, 5).Value = 18.014444444328547

$ws.Cells.Item(28, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(28, 2).Value = "402号直流"
$ws.Cells.Item(28, 3).Value = 45934.597407407404
$ws.Cells.Item(28, 4).Value = 45935.328888888886
$ws.Cells.Item(28, 5).Value = 17.555555555562023

$ws.Cells.Item(29, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(29, 2).Value = "306号直流"
$ws.Cells.Item(29, 3).Value = 45934.600995370369
$ws.Cells.Item(29, 4).Value = 45935.328888888886
$ws.Cells.Item(29, 5).Value = 17.469444444403052

$ws.Cells.Item(30, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(30, 2).Value = "905号直流"
$ws.Cells.Item(30, 3).Value = 45934.63821759259
$ws.Cells.Item(30, 4).Value = 45935.328888888886
$ws.Cells.Item(30, 5).Value = 16.57611111109145

$ws.Cells.Item(31, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(31, 2).Value = "A05号直流"
$ws.Cells.Item(31, 3).Value = 45934.645196759258
$ws.Cells.Item(31, 4).Value = 45935.328888888886
$ws.Cells.Item(31, 5).Value = 16.408611111051869

$ws.Cells.Item(32, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(32, 2).Value = "B04号直流"
$ws.Cells.Item(32, 3).Value = 45934.654016203705
$ws.Cells.Item(32, 4).Value = 45935.328888888886
$ws.Cells.Item(32, 5).Value = 16.19694444432389

$ws.Cells.Item(33, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(33, 2).Value = "212号直流"
$ws.Cells.Item(33, 3).Value = 45934.672256944446
$ws.Cells.Item(33, 4).Value = 45935.328888888886
$ws.Cells.Item(33, 5).Value = 15.759166666539386

$ws.Cells.Item(34, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(34, 2).Value = "206号直流"
$ws.Cells.Item(34, 3).Value = 45934.675659722219
$ws.Cells.Item(34, 4).Value = 45935.328888888886
$ws.Cells.Item(34, 5).Value = 15.677499999990687

$ws.Cells.Item(35, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(35, 2).Value = "105号直流"
$ws.Cells.Item(35, 3).Value = 45934.695162037038
$ws.Cells.Item(35, 4).Value = 45935.328888888886
$ws.Cells.Item(35, 5).Value = 15.209444444335531

$ws.Cells.Item(36, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(36, 2).Value = "110号直流"
$ws.Cells.Item(36, 3).Value = 45934.707303240742
$ws.Cells.Item(36, 4).Value = 45935.328888888886
$ws.Cells.Item(36, 5).Value = 14.918055555433966

$ws.Cells.Item(37, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(37, 2).Value = "108号直流"
$ws.Cells.Item(37, 3).Value = 45934.734479166669
$ws.Cells.Item(37, 4).Value = 45935.328888888886
$ws.Cells.Item(37, 5).Value = 14.265833333192859

$ws.Cells.Item(38, 1).Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Cells.Item(38, 2).Value = "804号直流"
$ws.Cells.Item(38, 3).Value = 45934.752916666665
$ws.Cells.Item(38, 4).Value = 45935.328888888886
$ws.Cells.Item(38, 5).Value = 13.823333333304618

$ws.Cells.Item(39, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(39, 2).Value = "401号直流"
$ws.Cells.Item(39, 3).Value = 45934.772418981483
$ws.Cells.Item(39, 4).Value = 45935.328888888886
$ws.Cells.Item(39, 5).Value = 13.355277777649462

$ws.Cells.Item(40, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(40, 2).Value = "202号直流"
$ws.Cells.Item(40, 3).Value = 45934.774108796293
$ws.Cells.Item(40, 4).Value = 45935.328888888886
$ws.Cells.Item(40, 5).Value = 13.314722222217824

$ws.Cells.Item(41, 1).Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Cells.Item(41, 2).Value = "011A号直流"
$ws.Cells.Item(41, 3).Value = 45934.798344907409
$ws.Cells.Item(41, 4).Value = 45935.328888888886
$ws.Cells.Item(41, 5).Value = 12.733055555436295

$ws.Cells.Item(42, 1).Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Cells.Item(42, 2).Value = "301号直流"
$ws.Cells.Item(42, 3).Value = 45934.814062500001
$ws.Cells.Item(42, 4).Value = 45935.328888888886
$ws.Cells.Item(42, 5).Value = 12.35583333321847

$ws.Cells.Item(43, 1).Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Cells.Item(43, 2).Value = "205号直流"
$ws.Cells.Item(43, 3).Value = 45934.823310185187
$ws.Cells.Item(43, 4).Value = 45935.328888888886
$ws.Cells.Item(43, 5).Value = 12.133888888754882

[void]$ws.Range("G11").Select()